# Connected Office Test Data.xlsx
# "Finished with Update and Delete tests"
#
# - Reset the "Create Test Passed" (B) and "Read Test Passed" (C) columns
#   on the "Test Results" sheet back to FALSE for every data row (2-24).
# - Move the remembered cell-selection on the "Device" sheet from D36 to C38
#   (without changing which sheet/tab is actually active - "Test Results"
#   stays the active tab, same as before the edit).

$wb = $excel.ActiveWorkbook

$testResults = $wb.Worksheets.Item("Test Results")
$device      = $wb.Worksheets.Item("Device")

# Reset Create/Read test-passed flags for rows 2 through 24.
$testResults.Range("B2:C24").Value = $false

# Update the stored selection on the Device sheet to C38, then restore
# "Test Results" as the active sheet so the active tab is unchanged.
$device.Activate()
$device.Range("C38").Select()
$testResults.Activate()
